$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: swap row2 / row3 identity ---
$overview.Range("A2").Value = "d1b96b37-2ce9-4e40-9935-38c6f7eeabb9.md"
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"

$overview.Range("A3").Value = "86ac0e1d-6bd0-4771-8bb6-2c05bf41ecdf.md"
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet ---
$zhcn.Range("A2").Value = "d1b96b37-2ce9-4e40-9935-38c6f7eeabb9.md"
$zhcn.Range("B2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C2").Value = "d1b96b37-2ce9-4e40-9935-38c6f7eeabb9.ba7f2da19cb425112cef85ba3e86e8ef8449c698.zh-cn.xlf"
$zhcn.Range("D2").Value = "2016-01-17 06:20:22"
$zhcn.Range("E2").Value = "d1b96b37-2ce9-4e40-9935-38c6f7eeabb9.md"
$zhcn.Range("F2").Value = "d1b96b37-2ce9-4e40-9935-38c6f7eeabb9.ba7f2da19cb425112cef85ba3e86e8ef8449c698.zh-cn.xlf"
$zhcn.Range("G2").Value = "2016-01-17 06:21:12"
$zhcn.Range("H2").Value = "Include"

$zhcn.Range("A3").Value = "86ac0e1d-6bd0-4771-8bb6-2c05bf41ecdf.md"
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "86ac0e1d-6bd0-4771-8bb6-2c05bf41ecdf.2cd4092bdafd04eeea090e000fd29a7c6105d867.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-01-17 06:22:14"
$zhcn.Range("E3").Value = "86ac0e1d-6bd0-4771-8bb6-2c05bf41ecdf.md"
$zhcn.Range("F3").Value = "86ac0e1d-6bd0-4771-8bb6-2c05bf41ecdf.2cd4092bdafd04eeea090e000fd29a7c6105d867.zh-cn.xlf"
$zhcn.Range("G3").Value = "2016-01-17 06:21:12"
$zhcn.Range("H3").Value = "Include"

# --- de-de sheet ---
$dede.Range("A2").Value = "d1b96b37-2ce9-4e40-9935-38c6f7eeabb9.md"
$dede.Range("B2").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "d1b96b37-2ce9-4e40-9935-38c6f7eeabb9.ba7f2da19cb425112cef85ba3e86e8ef8449c698.de-de.xlf"
$dede.Range("D2").Value = "2016-01-17 06:20:34"
$dede.Range("E2").Value = "d1b96b37-2ce9-4e40-9935-38c6f7eeabb9.md"
$dede.Range("F2").Value = "d1b96b37-2ce9-4e40-9935-38c6f7eeabb9.ba7f2da19cb425112cef85ba3e86e8ef8449c698.de-de.xlf"
$dede.Range("G2").Value = "2016-01-17 06:21:30"
$dede.Range("H2").Value = "Include"

$dede.Range("A3").Value = "86ac0e1d-6bd0-4771-8bb6-2c05bf41ecdf.md"
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "86ac0e1d-6bd0-4771-8bb6-2c05bf41ecdf.2cd4092bdafd04eeea090e000fd29a7c6105d867.de-de.xlf"
$dede.Range("D3").Value = "2016-01-17 06:22:24"
$dede.Range("E3").Value = "86ac0e1d-6bd0-4771-8bb6-2c05bf41ecdf.md"
$dede.Range("F3").Value = "86ac0e1d-6bd0-4771-8bb6-2c05bf41ecdf.2cd4092bdafd04eeea090e000fd29a7c6105d867.de-de.xlf"
$dede.Range("G3").Value = "2016-01-17 06:21:30"
$dede.Range("H3").Value = "Include"

# --- Update hyperlinks display text (and keep same relationship targets) ---
foreach ($ws in @($overview, $zhcn, $dede)) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address -eq "$($hl.Range.Worksheet.Range('A2').Address)") { }
    }
}
